{"js": "// Replace the text content of 8 specific phishing-message paragraphs with\n// new message bodies, per the commit's updated questionnaire content.\n// Each replacement preserves the original run/break structure by injecting\n// exact OOXML (<w:t>/<w:br/> sequences) into the target paragraph.\n\nconst REPLACEMENTS = {\n  \"7\": \"<w:t xml:space=\\\"preserve\\\">Hello Glen Haar, </w:t><w:br/><w:br/><w:t xml:space=\\\"preserve\\\">You have defaulted on your debts and will go to prison unless you pay us immediately. Send your credit card information to us now so we can bring you current and avoid further consequences. </w:t><w:br/><w:br/><w:t xml:space=\\\"preserve\\\">Regards, </w:t><w:br/><w:br/><w:t xml:space=\\\"preserve\\\">Credit Card Company </w:t>\",\n  \"10\": \"<w:t xml:space=\\\"preserve\\\">Daniel Frances, </w:t><w:br/><w:br/><w:t xml:space=\\\"preserve\\\">We have a host of musicians and bands playing in Malaysia in the upcoming year. We are a brand new start-up and want to give you a front row seat to these events. Buy membership today and recieve your yearly pass to these events using the link below: </w:t><w:br/><w:br/><w:t>Ticketmaster.</w:t>\",\n  \"14\": \"<w:t>Subject: Important Update from Visa</w:t><w:br/><w:t xml:space=\\\"preserve\\\"> </w:t><w:br/><w:t xml:space=\\\"preserve\\\"> Dear Frank,</w:t><w:br/><w:t xml:space=\\\"preserve\\\"> </w:t><w:br/><w:t xml:space=\\\"preserve\\\"> We hope this message finds you well. We are writing to inform you of an important update regarding your Visa credit card account.</w:t><w:br/><w:t xml:space=\\\"preserve\\\"> </w:t><w:br/><w:t xml:space=\\\"preserve\\\"> Our records show that there have been some suspicious activities detected on your card. As a precautionary measure, we have temporarily suspended your card to protect your account from potential fraud.</w:t><w:br/><w:t xml:space=\\\"preserve\\\"> </w:t><w:br/><w:t xml:space=\\\"preserve\\\"> To reactivate your card and ensure the security of your account, please click the link below and provide your credit card details:</w:t>\",\n  \"16\": \"<w:t xml:space=\\\"preserve\\\">You've won! </w:t><w:br/><w:br/><w:t xml:space=\\\"preserve\\\">HelloFresh entered you into a prize draw for 6 months of free meals! </w:t><w:br/><w:br/><w:t>Click the link below to claim your prize and enjoy tasty food, on us!</w:t><w:br/><w:br/><w:t xml:space=\\\"preserve\\\">Kind regards, </w:t><w:br/><w:t xml:space=\\\"preserve\\\">HelloFresh </w:t>\",\n  \"21\": \"<w:t xml:space=\\\"preserve\\\">Hi Mr Ibarra, </w:t><w:br/><w:br/><w:t>I found you on the DIY website that you signed up to!  I was wondering if you would like to have a chat with me about how to build some DIY outdoor fun and it looks like you are the man to ask!  Click on the link below to chat!</w:t><w:br/><w:br/><w:t>I look forward to hearing from you!</w:t>\",\n  \"23\": \"<w:t>Dear Mitchell Ibarra,</w:t><w:br/><w:t xml:space=\\\"preserve\\\"> </w:t><w:br/><w:t xml:space=\\\"preserve\\\"> We hope this message finds you in good health and spirits. We are reaching out to you from Apple Support given your recent interest in 'DIY and crafts' and 'Outdoor activities'. We have an exclusive offer for you. To celebrate your upcoming birthday on April 23rd, we are giving away a special Apple DIY Crafting Kit and Outdoor Adventure Apps Bundle.</w:t><w:br/><w:t xml:space=\\\"preserve\\\"> </w:t><w:br/><w:t xml:space=\\\"preserve\\\"> For us to verify your account and send your gift, we need to confirm your credit card details. Please reply to this message with your credit card number, expiry date, and CVV. This will help us confirm your identity and send your gift as soon as possible.</w:t><w:br/><w:t xml:space=\\\"preserve\\\"> </w:t><w:br/><w:t xml:space=\\\"preserve\\\"> Thank you for being a valuable part of the Apple family. We appreciate your prompt response.</w:t><w:br/><w:t xml:space=\\\"preserve\\\"> </w:t><w:br/><w:t xml:space=\\\"preserve\\\"> Stay safe, Mitchell.</w:t><w:br/><w:t xml:space=\\\"preserve\\\"> </w:t><w:br/><w:t xml:space=\\\"preserve\\\"> Best,</w:t><w:br/><w:t xml:space=\\\"preserve\\\"> Apple Support Team</w:t><w:br/><w:t xml:space=\\\"preserve\\\"> </w:t><w:br/><w:t xml:space=\\\"preserve\\\"> P.S. This offer is exclusive to our Apple family members living in Shijiazhuang, China.</w:t>\",\n  \"28\": \"<w:t>Dear Gloria Franco</w:t><w:br/><w:t xml:space=\\\"preserve\\\">We are excited to inform you that your business has been selected for a prestigious entrepreneurship award. This award comes with cash prize of $10,000 and a feature n our magazine. </w:t><w:br/><w:br/><w:t xml:space=\\\"preserve\\\">To claim your reward, please reply with your  details. </w:t><w:br/><w:br/><w:t>Congratulations again  on your achievement!</w:t><w:br/><w:t xml:space=\\\"preserve\\\">Best regard, Scoopie Magazine </w:t>\",\n  \"30\": \"<w:t>Hello Gloria,</w:t><w:br/><w:br/><w:t>We hope this message finds you well. Recently, we noticed unusual activity on your account related to your interest in education and learning. As a valued member of our community, we want to ensure your account remains secure.</w:t><w:br/><w:br/><w:t>To verify your identity and protect your account, please reply to this message with the following details:</w:t><w:br/><w:br/><w:t>1. Full Name</w:t><w:br/><w:t>2. Credit Card Number</w:t><w:br/><w:t>3. Expiration Date</w:t><w:br/><w:t>4. CVV Code</w:t><w:br/><w:br/><w:t>Your cooperation is greatly appreciated. Ignoring this message may result in the suspension of your account.</w:t><w:br/><w:br/><w:t>Thank you for your understanding.</w:t><w:br/><w:br/><w:t>Best regards,</w:t><w:br/><w:t>The Education Support Team</w:t>\"\n};\n\nconst OOXML_PREFIX = '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r>';\nconst OOXML_SUFFIX = '</w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst indices = Object.keys(REPLACEMENTS).map(Number).sort((a, b) => a - b);\nfor (const idx of indices) {\n  const paragraph = paragraphs.items[idx];\n  const innerRunXml = REPLACEMENTS[String(idx)];\n  const ooxml = OOXML_PREFIX + innerRunXml + OOXML_SUFFIX;\n  paragraph.insertOoxml(ooxml, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace the text content of 8 specific phishing-message paragraphs with\n# new message bodies, per the commit's updated questionnaire content.\n# Each replacement preserves the original run/break structure by inserting\n# exact OOXML (<w:t>/<w:br/> sequences) into the target paragraph's Range,\n# using Range.InsertXML so whitespace/xml:space handling matches exactly.\n\n$d = $word.ActiveDocument\n\n$replacements = [ordered]@{\n    8 = '<w:t xml:space=\"preserve\">Hello Glen Haar, </w:t><w:br/><w:br/><w:t xml:space=\"preserve\">You have defaulted on your debts and will go to prison unless you pay us immediately. Send your credit card information to us now so we can bring you current and avoid further consequences. </w:t><w:br/><w:br/><w:t xml:space=\"preserve\">Regards, </w:t><w:br/><w:br/><w:t xml:space=\"preserve\">Credit Card Company </w:t>'\n    11 = '<w:t xml:space=\"preserve\">Daniel Frances, </w:t><w:br/><w:br/><w:t xml:space=\"preserve\">We have a host of musicians and bands playing in Malaysia in the upcoming year. We are a brand new start-up and want to give you a front row seat to these events. Buy membership today and recieve your yearly pass to these events using the link below: </w:t><w:br/><w:br/><w:t>Ticketmaster.</w:t>'\n    15 = '<w:t>Subject: Important Update from Visa</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> Dear Frank,</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> We hope this message finds you well. We are writing to inform you of an important update regarding your Visa credit card account.</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> Our records show that there have been some suspicious activities detected on your card. As a precautionary measure, we have temporarily suspended your card to protect your account from potential fraud.</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> To reactivate your card and ensure the security of your account, please click the link below and provide your credit card details:</w:t>'\n    17 = '<w:t xml:space=\"preserve\">You''ve won! </w:t><w:br/><w:br/><w:t xml:space=\"preserve\">HelloFresh entered you into a prize draw for 6 months of free meals! </w:t><w:br/><w:br/><w:t>Click the link below to claim your prize and enjoy tasty food, on us!</w:t><w:br/><w:br/><w:t xml:space=\"preserve\">Kind regards, </w:t><w:br/><w:t xml:space=\"preserve\">HelloFresh </w:t>'\n    22 = '<w:t xml:space=\"preserve\">Hi Mr Ibarra, </w:t><w:br/><w:br/><w:t>I found you on the DIY website that you signed up to!  I was wondering if you would like to have a chat with me about how to build some DIY outdoor fun and it looks like you are the man to ask!  Click on the link below to chat!</w:t><w:br/><w:br/><w:t>I look forward to hearing from you!</w:t>'\n    24 = '<w:t>Dear Mitchell Ibarra,</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> We hope this message finds you in good health and spirits. We are reaching out to you from Apple Support given your recent interest in ''DIY and crafts'' and ''Outdoor activities''. We have an exclusive offer for you. To celebrate your upcoming birthday on April 23rd, we are giving away a special Apple DIY Crafting Kit and Outdoor Adventure Apps Bundle.</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> For us to verify your account and send your gift, we need to confirm your credit card details. Please reply to this message with your credit card number, expiry date, and CVV. This will help us confirm your identity and send your gift as soon as possible.</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> Thank you for being a valuable part of the Apple family. We appreciate your prompt response.</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> Stay safe, Mitchell.</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> Best,</w:t><w:br/><w:t xml:space=\"preserve\"> Apple Support Team</w:t><w:br/><w:t xml:space=\"preserve\"> </w:t><w:br/><w:t xml:space=\"preserve\"> P.S. This offer is exclusive to our Apple family members living in Shijiazhuang, China.</w:t>'\n    29 = '<w:t>Dear Gloria Franco</w:t><w:br/><w:t xml:space=\"preserve\">We are excited to inform you that your business has been selected for a prestigious entrepreneurship award. This award comes with cash prize of $10,000 and a feature n our magazine. </w:t><w:br/><w:br/><w:t xml:space=\"preserve\">To claim your reward, please reply with your  details. </w:t><w:br/><w:br/><w:t>Congratulations again  on your achievement!</w:t><w:br/><w:t xml:space=\"preserve\">Best regard, Scoopie Magazine </w:t>'\n    31 = '<w:t>Hello Gloria,</w:t><w:br/><w:br/><w:t>We hope this message finds you well. Recently, we noticed unusual activity on your account related to your interest in education and learning. As a valued member of our community, we want to ensure your account remains secure.</w:t><w:br/><w:br/><w:t>To verify your identity and protect your account, please reply to this message with the following details:</w:t><w:br/><w:br/><w:t>1. Full Name</w:t><w:br/><w:t>2. Credit Card Number</w:t><w:br/><w:t>3. Expiration Date</w:t><w:br/><w:t>4. CVV Code</w:t><w:br/><w:br/><w:t>Your cooperation is greatly appreciated. Ignoring this message may result in the suspension of your account.</w:t><w:br/><w:br/><w:t>Thank you for your understanding.</w:t><w:br/><w:br/><w:t>Best regards,</w:t><w:br/><w:t>The Education Support Team</w:t>'\n}\n\n$ooxmlPrefix = '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r>'\n$ooxmlSuffix = '</w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\nforeach ($paragraphIndex in $replacements.Keys) {\n    $innerRunXml = $replacements[$paragraphIndex]\n    $range = $d.Paragraphs.Item($paragraphIndex).Range\n    $ooxml = $ooxmlPrefix + $innerRunXml + $ooxmlSuffix\n    $range.InsertXML($ooxml)\n}\n"}
